$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.23"
$ws.Range("E2").Value = "'-1.24%"
$ws.Range("D3").Value = "'35.83"
$ws.Range("E3").Value = "'-5.05%"
$ws.Range("D4").Value = "'5.113"
$ws.Range("E4").Value = "'-0.69%"
$ws.Range("D5").Value = "'0.07693"
$ws.Range("E5").Value = "'-2.77%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.391"
$ws.Range("E6").Value = "'-0.53%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.296"
$ws.Range("E7").Value = "'0.32%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.840"
$ws.Range("E8").Value = "'-3.19%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.952"
$ws.Range("E9").Value = "'-4.70%"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9207"
$ws.Range("E10").Value = "'-0.10%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1114"
$ws.Range("E11").Value = "'-7.75%"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1857"
$ws.Range("E12").Value = "'-4.08%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.08755"
$ws.Range("E13").Value = "'-5.14%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03337"
$ws.Range("E14").Value = "'1.03%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09519"
$ws.Range("E15").Value = "'-1.17%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001394"
$ws.Range("E16").Value = "'0.97%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006225"
$ws.Range("E17").Value = "'6.35%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.365"
$ws.Range("E18").Value = "'-4.34%"
$ws.Range("D19").Value = "'0.3444"
$ws.Range("E19").Value = "'1.26%"
$ws.Range("D20").Value = "'6.304"
$ws.Range("E20").Value = "'19.35%"
$ws.Range("D21").Value = "'0.1292"
$ws.Range("E21").Value = "'1.52%"
$ws.Range("E22").Value = "'-10.66%"
$ws.Range("D23").Value = "'0.04347"
$ws.Range("E23").Value = "'-0.28%"
$ws.Range("E24").Value = "'-3.51%"
$ws.Range("D25").Value = "'0.004253"
$ws.Range("E25").Value = "'-1.40%"
$ws.Range("D26").Value = "'0.0001332"
$ws.Range("E26").Value = "'9.20%"
$ws.Range("D27").Value = "'0.0002905"
$ws.Range("D39").Value = "'0.02093"
$ws.Range("E39").Value = "'-1.68%"
$ws.Range("D40").Value = "'0.04914"
$ws.Range("E40").Value = "'-4.99%"
$ws.Range("D41").Value = "'0.007541"
$ws.Range("E41").Value = "'-1.02%"
$ws.Range("E42").Value = "'-1.18%"
$ws.Range("D43").Value = "'0.008576"
$ws.Range("E43").Value = "'-5.62%"
$ws.Range("D44").Value = "'0.002073"
$ws.Range("E44").Value = "'3.16%"
$ws.Range("D45").Value = "'0.008397"
$ws.Range("E45").Value = "'-2.36%"
$ws.Range("D46").Value = "'0.00006467"
$ws.Range("E46").Value = "'-3.30%"
$ws.Range("E47").Value = "'0.21%"
$ws.Range("D48").Value = "'0.003301"
$ws.Range("E48").Value = "'14.67%"
$ws.Range("D49").Value = "'0.001445"
$ws.Range("E49").Value = "'20.51%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.21%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.21%"

# Restore default (unstyled) formatting for text-forced numeric/percentage cells
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Style = "Normal"

